$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the header row: "<label>_old" -> "<label>_FV2410", "<label>_new" -> "<label>_FV2504" ---
$labels = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = "$($labels[$i])_FV2410"
}
# Column 11 ("diff") is untouched.
for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = "$($labels[$i])_FV2504"
}

# --- Turn the data range into an Excel Table (ListObject) with autofilter ---
$rng = $ws.Range("A1:U76")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

# --- Freeze the header row (split below row 1) ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
